# test_metafile_2boxes.xlsx commit:
#   "Updated by adding log / also changed compression to earlier form"
#
# Semantic changes applied:
#   1. Rename "test_2boxes" -> "test_metafile_2boxes" and
#      "test_1box" -> "test_metafile_1box".
#   2. Shorten the duplicated "sig-"/"ttl-" header labels (K,L,N,O,P on
#      row 1) on every data sheet - "sig-blue"->"blue", "sig-uv"->"uv",
#      "ttl-licks"->"licks", "ttl-distractors"->"distractors",
#      "ttl-distracted"->"distracted". ("tick" in column M is untouched.)
#   3. Update the saved view/selection state: the big metafile sheet is
#      scrolled down and rows 88:89 are selected; the (renamed) 2-box
#      sheet becomes the active tab with K1 selected; the (renamed)
#      1-box sheet ends with F1 selected.

$wb = $excel.ActiveWorkbook

$wsMeta  = $wb.Worksheets.Item("thph12_metafile")
$ws2box  = $wb.Worksheets.Item("test_2boxes")
$ws1box  = $wb.Worksheets.Item("test_1box")

# --- 1. Rename sheets ------------------------------------------------
$ws2box.Name = "test_metafile_2boxes"
$ws1box.Name = "test_metafile_1box"

# --- 2. Fix up the repeated header row on every data sheet -----------
foreach ($ws in @($wsMeta, $ws2box, $ws1box)) {
    $ws.Range("K1").Value = "blue"
    $ws.Range("L1").Value = "uv"
    $ws.Range("N1").Value = "licks"
    $ws.Range("O1").Value = "distractors"
    $ws.Range("P1").Value = "distracted"
}

# --- 3. Restore view / selection state --------------------------------
# Big metafile sheet: scrolled so row 64 is at the top, rows 88:89
# selected, no longer the active tab.
[void]$wsMeta.Activate()
$excel.ActiveWindow.ScrollRow = 64
$excel.ActiveWindow.ScrollColumn = 1
[void]$wsMeta.Range("A88:XFD89").Select()

# 1-box sheet: single cell F1 selected.
[void]$ws1box.Activate()
[void]$ws1box.Range("F1").Select()

# 2-box sheet: single cell K1 selected, and this is the sheet left
# active/selected when the workbook is saved.
[void]$ws2box.Activate()
[void]$ws2box.Range("K1").Select()
